$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 989
$ws.Range("E2").Value = 21
$ws.Range("F2").Value = 21
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 6
$ws.Range("I2").Value = 6
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 1010
$ws.Range("L2").Value = 603
$ws.Range("M2").Value = 407
$ws.Range("N2").Value = 407
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 143
$ws.Range("Q2").Value = 97
$ws.Range("R2").Value = -24
$ws.Range("S2").Value = -15
$ws.Range("T2").Value = 15
$ws.Range("U2").Value = 82
$ws.Range("V2").Value = 227
$ws.Range("W2").Value = 2.08
$ws.Range("X2").Value = 0.6
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 0.59
$ws.Range("AA2").Value = 148.35
$ws.Range("AB2").Value = 145.45
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 91.12
$ws.Range("AE2").Value = 1423
$ws.Range("AF2").Value = 1.34
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 28595074

# Row 3
$ws.Range("D3").Value = 965
$ws.Range("E3").Value = -78
$ws.Range("F3").Value = -78
$ws.Range("G3").Value = -96
$ws.Range("H3").Value = -109
$ws.Range("I3").Value = -108
$ws.Range("J3").Value = -1
$ws.Range("K3").Value = 894
$ws.Range("L3").Value = 586
$ws.Range("M3").Value = 308
$ws.Range("N3").Value = 308
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 148
$ws.Range("Q3").Value = 109
$ws.Range("R3").Value = -40
$ws.Range("S3").Value = 6
$ws.Range("T3").Value = 35
$ws.Range("U3").Value = 74
$ws.Range("V3").Value = 219
$ws.Range("W3").Value = -8.109999999999999
$ws.Range("X3").Value = -11.26
$ws.Range("Y3").Value = -30.21
$ws.Range("Z3").Value = -11.42
$ws.Range("AA3").Value = 190.4
$ws.Range("AB3").Value = 108.48
$ws.Range("AC3").Value = -366
$ws.Range("AD3").Value = -8.08
$ws.Range("AE3").Value = 1041
$ws.Range("AF3").Value = 2.84
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 29558465

# Row 4
$ws.Range("D4").Value = 971
$ws.Range("E4").Value = -35
$ws.Range("F4").Value = -35
$ws.Range("G4").Value = -58
$ws.Range("H4").Value = -77
$ws.Range("I4").Value = -76
$ws.Range("J4").Value = -1
$ws.Range("K4").Value = 769
$ws.Range("L4").Value = 523
$ws.Range("M4").Value = 246
$ws.Range("N4").Value = 246
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 155
$ws.Range("Q4").Value = -5
$ws.Range("R4").Value = -46
$ws.Range("S4").Value = -8
$ws.Range("T4").Value = 48
$ws.Range("U4").Value = -52
$ws.Range("V4").Value = 195
$ws.Range("W4").Value = -3.6
$ws.Range("X4").Value = -7.94
$ws.Range("Y4").Value = -27.61
$ws.Range("Z4").Value = -9.27
$ws.Range("AA4").Value = 212.26
$ws.Range("AB4").Value = 59.38
$ws.Range("AC4").Value = -250
$ws.Range("AD4").Value = -6.49
$ws.Range("AE4").Value = 797
$ws.Range("AF4").Value = 2.04
$ws.Range("AG4").ClearContents()
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 30907210

# Row 5
$ws.Range("D5").Value = 1141
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = -4
$ws.Range("H5").Value = -6
$ws.Range("I5").Value = -8
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 791
$ws.Range("L5").Value = 554
$ws.Range("M5").Value = 237
$ws.Range("N5").Value = 235
$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 155
$ws.Range("Q5").Value = -13
$ws.Range("R5").Value = -63
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 52
$ws.Range("U5").Value = -65
$ws.Range("V5").Value = 200
$ws.Range("W5").Value = 0.58
$ws.Range("X5").Value = -0.49
$ws.Range("Y5").Value = -3.12
$ws.Range("Z5").Value = -0.71
$ws.Range("AA5").Value = 234.14
$ws.Range("AB5").Value = 51.98
$ws.Range("AC5").Value = -24
$ws.Range("AD5").Value = -45.7
$ws.Range("AE5").Value = 760
$ws.Range("AF5").Value = 1.46
$ws.Range("AG5").ClearContents()
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 30907210

# Row 6
$ws.Range("D6").Value = 1086
$ws.Range("E6").Value = 12
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = -1
$ws.Range("I6").Value = -2
$ws.Range("J6").ClearContents()
$ws.Range("K6").Value = 785
$ws.Range("L6").Value = 551
$ws.Range("M6").Value = 234
$ws.Range("N6").Value = 232
$ws.Range("O6").ClearContents()
$ws.Range("P6").Value = 155
$ws.Range("Q6").Value = 29
$ws.Range("R6").Value = -17
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 16
$ws.Range("U6").Value = 13
$ws.Range("V6").Value = 199
$ws.Range("W6").Value = 1.14
$ws.Range("X6").Value = -0.12
$ws.Range("Y6").Value = -0.65
$ws.Range("Z6").Value = -0.16
$ws.Range("AA6").Value = 234.95
$ws.Range("AB6").Value = 50.69
$ws.Range("AC6").Value = -5
$ws.Range("AD6").Value = -273.2
$ws.Range("AE6").Value = 751
$ws.Range("AF6").Value = 1.8
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").ClearContents()
$ws.Range("AJ6").Value = 30907210

# Row 7: clear all data columns
$ws.Range("D7:AJ7").ClearContents()

# Row 8: clear all data columns
$ws.Range("D8:AJ8").ClearContents()

# Row 9: clear all data columns
$ws.Range("D9:AJ9").ClearContents()

Write-Host "edits applied"